$wb = $excel.ActiveWorkbook

# 1) Remove the stale duplicate outage row (R4 / JED0124) from sheet "R1".
#    It duplicated row 5's data with a different Hub Site value, and is no
#    longer present in the refreshed export.
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Rows.Item(6).Delete()

# 2) Refresh the "Elapsed Duration(Hrs)" values (column G) across sheets to
#    reflect the later report-generation timestamp.
$ws1.Range("G2").Value = "3946:53:15"
$ws1.Range("G3").Value = "86:25:53"
$ws1.Range("G4").Value = "109:25:53"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12128:16:53"
$ws2.Range("G3").Value = "3258:00:22"
$ws2.Range("G4").Value = "496:11:56"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2974:06:42"
$ws4.Range("G3").Value = "201:18:57"
$ws4.Range("G4").Value = "89:31:22"
$ws4.Range("G5").Value = "87:08:55"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "448:05:41"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "88:37:59"
